# Master_CON_GS_CO2_DK.xlsx — Gasera CO2 -> C-CO2 transformation
# (CCH4, NN2O and CCO2 Gasera transformations in Master_GHG_2023)
#
# Renames the Gasera CO2 flux headers to the "C-CO2" naming and rescales
# the corresponding flux values (columns C & D) by the CO2 -> C-CO2
# conversion factor (44/12 * 9/11 = 121/9 = 13.444444444444445), leaving
# the chromatograph column (E) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames -------------------------------------------------
$ws.Range("C1").Value = "avg_Gasera_CCO2_flux_mgm2h"
$ws.Range("D1").Value = "avg_Gasera_CCO2_flux_mgm2h_cor"

# --- Rescale the Gasera CO2 flux values (rows 2-16) ------------------
$factor = 13.444444444444445

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $rawValue = $ws.Cells.Item($row, 3).Value2
    if ($rawValue -ne $null) {
        $converted = $rawValue / $factor
        $ws.Cells.Item($row, 3).Value = $converted
        $ws.Cells.Item($row, 4).Value = $converted
    }
}
